$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row with two new predicted-factor columns (D, E),
# copying the bold/bordered header style already used by C1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Refresh column C values (rows 2-6) with the recalculated predictions
$ws.Range("C2").Value = -4.697124697347145
$ws.Range("C3").Value = -1.172924015787017
$ws.Range("C4").Value = -0.07027665786814449
$ws.Range("C5").Value = -0.4595607842740025
$ws.Range("C6").Value = -0.1300780636132118

# Fill in the new column D values (rows 2-6)
$ws.Range("D2").Value = -4.334917804198449
$ws.Range("D3").Value = -1.133843532881532
$ws.Range("D4").Value = 0.00361369713132461
$ws.Range("D5").Value = -0.2015390863418249
$ws.Range("D6").Value = -0.1323400936722163

# Fill in the new column E values (rows 2-6)
$ws.Range("E2").Value = -3.972531378745367
$ws.Range("E3").Value = -1.081707472675687
$ws.Range("E4").Value = 0.0592240151895012
$ws.Range("E5").Value = 0.00001166751880576564
$ws.Range("E6").Value = -0.131291933233199
